$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on the Price (D) column cells being updated so that
# numeric-looking strings (e.g. "1.000", "0.9997") are stored as text,
# exactly as they appear in the source data, instead of being parsed as numbers.
$priceCells = @(
    "D2", "D3", "D4", "D5", "D6", "D7", "D8", "D9", "D10", "D11",
    "D12", "D13", "D14", "D15", "D16", "D17", "D18", "D19", "D20", "D21",
    "D22", "D23", "D24", "D25", "D26", "D27", "D28", "D29", "D30", "D31",
    "D32", "D33", "D34", "D35", "D36", "D37", "D38", "D39", "D40", "D41",
    "D43", "D44", "D45", "D46", "D47", "D48", "D49"
)
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# --- Apply cell value updates row by row (matches the source diff) ---

# Row 2
$ws.Range("D2").Value = "30.040.38"
$ws.Range("E2").Value = "  -1.56%  "

# Row 3
$ws.Range("D3").Value = "1.855.29"
$ws.Range("E3").Value = "  +0.10%  "

# Row 4
$ws.Range("D4").Value = "0.9997"
$ws.Range("E4").Value = "  -0.07%  "

# Row 5
$ws.Range("D5").Value = "235.60"
$ws.Range("E5").Value = "  +0.70%  "

# Row 6
$ws.Range("D6").Value = "1.000"

# Row 7
$ws.Range("D7").Value = "0.4638"
$ws.Range("E7").Value = "  -1.43%  "

# Row 8
$ws.Range("D8").Value = "0.2766"
$ws.Range("E8").Value = "  +0.55%  "

# Row 9
$ws.Range("D9").Value = "0.06406"
$ws.Range("E9").Value = "  +1.09%  "

# Row 10
$ws.Range("D10").Value = "18.33"
$ws.Range("E10").Value = "  +3.74%  "

# Row 11
$ws.Range("D11").Value = "97.35"
$ws.Range("E11").Value = "  +14.93%  "

# Row 12
$ws.Range("D12").Value = "1.848.53"
$ws.Range("E12").Value = "  -0.19%  "

# Row 13
$ws.Range("D13").Value = "0.07521"
$ws.Range("E13").Value = "  +1.49%  "

# Row 14
$ws.Range("D14").Value = "4.974"
$ws.Range("E14").Value = "  -1.05%  "

# Row 15
$ws.Range("D15").Value = "0.6225"
$ws.Range("E15").Value = "  -0.96%  "

# Row 16
$ws.Range("D16").Value = "295.47"
$ws.Range("E16").Value = "  +22.20%  "

# Row 17
$ws.Range("D17").Value = "29.975.57"
$ws.Range("E17").Value = "  -1.68%  "

# Row 18
$ws.Range("D18").Value = "1.000"
$ws.Range("E18").Value = "  -0.03%  "

# Row 19
$ws.Range("B19").Value = "ShibaInu"
$ws.Range("C19").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D19").Value = "0.000007370"
$ws.Range("E19").Value = "  +0.12%  "

# Row 20
$ws.Range("B20").Value = "Avalanche"
$ws.Range("C20").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D20").Value = "12.56"
$ws.Range("E20").Value = "  -1.09%  "

# Row 21
$ws.Range("D21").Value = "2.088.04"
$ws.Range("E21").Value = "  +0.05%  "

# Row 22
$ws.Range("D22").Value = "0.9992"
$ws.Range("E22").Value = "  +0.07%  "

# Row 23
$ws.Range("D23").Value = "4.992"
$ws.Range("E23").Value = "  +0.51%  "

# Row 24
$ws.Range("D24").Value = "6.070"
$ws.Range("E24").Value = "  +1.38%  "

# Row 25
$ws.Range("D25").Value = "164.54"
$ws.Range("E25").Value = "  +1.44%  "

# Row 26
$ws.Range("D26").Value = "9.065"
$ws.Range("E26").Value = "  -2.16%  "

# Row 27
$ws.Range("D27").Value = "19.14"
$ws.Range("E27").Value = "  +6.07%  "

# Row 28
$ws.Range("D28").Value = "1.927"
$ws.Range("E28").Value = "  +1.88%  "

# Row 29
$ws.Range("D29").Value = "0.1075"
$ws.Range("E29").Value = "  +5.91%  "

# Row 30
$ws.Range("D30").Value = "1.323"
$ws.Range("E30").Value = "  -3.17%  "

# Row 31
$ws.Range("D31").Value = "3.978"
$ws.Range("E31").Value = "  -1.49%  "

# Row 32
$ws.Range("D32").Value = "3.822"
$ws.Range("E32").Value = "  -0.90%  "

# Row 33
$ws.Range("D33").Value = "0.04857"
$ws.Range("E33").Value = "  -0.68%  "

# Row 34
$ws.Range("D34").Value = "0.7321"
$ws.Range("E34").Value = "  +3.41%  "

# Row 35
$ws.Range("D35").Value = "1.111"
$ws.Range("E35").Value = "  -2.59%  "

# Row 36
$ws.Range("D36").Value = "2.728"
$ws.Range("E36").Value = "  +0.86%  "

# Row 37
$ws.Range("D37").Value = "0.01891"
$ws.Range("E37").Value = "  -0.85%  "

# Row 38
$ws.Range("D38").Value = "2.643"
$ws.Range("E38").Value = "  -1.69%  "

# Row 39
$ws.Range("B39").Value = "Quant"
$ws.Range("C39").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D39").Value = "106.64"
$ws.Range("E39").Value = "  +1.21%  "

# Row 40
$ws.Range("B40").Value = "RenderToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D40").Value = "1.957"
$ws.Range("E40").Value = "  -1.30%  "

# Row 41
$ws.Range("B41").Value = "TrustWalletToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D41").Value = "0.8533"
$ws.Range("E41").Value = "  -2.61%  "

# Row 42
$ws.Range("E42").Value = "  +0.13%  "

# Row 43
$ws.Range("D43").Value = "5.725"
$ws.Range("E43").Value = "  +3.77%  "

# Row 44
$ws.Range("D44").Value = "0.4032"
$ws.Range("E44").Value = "  -1.38%  "

# Row 45
$ws.Range("D45").Value = "65.59"
$ws.Range("E45").Value = "  +4.58%  "

# Row 46
$ws.Range("D46").Value = "7.057"
$ws.Range("E46").Value = "  -2.54%  "

# Row 47
$ws.Range("D47").Value = "8.972"
$ws.Range("E47").Value = "  +4.41%  "

# Row 48
$ws.Range("D48").Value = "0.1188"
$ws.Range("E48").Value = "  -1.50%  "

# Row 49
$ws.Range("D49").Value = "33.73"
$ws.Range("E49").Value = "  +1.05%  "

# Row 50
$ws.Range("E50").Value = "  +0.19%  "

# Row 51
$ws.Range("E51").Value = "  +0.14%  "

# Restore default ("Normal") style on the Price cells so only the content
# differs from the original workbook (formatting/style stays the same).
foreach ($addr in $priceCells) {
    $ws.Range($addr).Style = "Normal"
}
